$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "biegen",
    "sterben",
    "spinnen",
    "machen",
    "feiern",
    "merken",
    "reiten",
    "arten",
    "sorgen",
    "zielen",
    "fahren",
    "weichen",
    "locken",
    "klingeln",
    "stammen",
    "schützen",
    "lügen",
    "schulden",
    "räumen",
    "lächeln",
    "planen",
    "klingen",
    "schwören",
    "treiben",
    "deuten",
    "landen",
    "fallen",
    "danken",
    "weinen",
    "scheitern",
    "wehtun",
    "mauern",
    "wollen",
    "werfen",
    "folgen",
    "bauen",
    "spüren",
    "dringen",
    "jagen",
    "scheiden",
    "rasen",
    "plaudern",
    "greifen",
    "saufen",
    "laufen",
    "schaden",
    "dienen",
    "schlucken",
    "gelten",
    "treffen",
    "feuern",
    "führen",
    "enden",
    "grüßen",
    "schütteln",
    "füllen",
    "wachsen",
    "schrecken",
    "zünden",
    "töten",
    "zählen",
    "rühren",
    "fassen",
    "sitzen",
    "jubeln",
    "wandern",
    "achten",
    "fließen",
    "irren",
    "wüten",
    "starten",
    "zögern",
    "malen",
    "sprengen",
    "brauchen",
    "brauen",
    "trennen",
    "helfen",
    "warnen",
    "seufzen",
    "rufen",
    "bluten",
    "kämpfen",
    "fischen",
    "proben",
    "fangen",
    "meistern",
    "suchen",
    "fragen",
    "bellen",
    "pressen",
    "stehlen",
    "hören",
    "leiden",
    "runden",
    "heulen",
    "ärgern",
    "gnaden",
    "rechnen",
    "kosten",
    "wundern",
    "trauen",
    "klettern",
    "flüchten",
    "sagen",
    "ehren",
    "holen",
    "bergen",
    "platzen",
    "spielen",
    "haben",
    "schwächen",
    "backen",
    "retten",
    "fällen",
    "siegen",
    "rücken",
    "freuen",
    "stoppen",
    "knarren",
    "bitten",
    "tauchen",
    "zahlen",
    "lesen",
    "dauern",
    "boxen",
    "warten",
    "stecken",
    "äußern",
    "schmecken",
    "altern",
    "liefern",
    "eignen",
    "füttern",
    "schreiten",
    "schenken",
    "flehen",
    "schwimmen",
    "wetten",
    "sperren",
    "münzen",
    "geben",
    "betteln",
    "schämen",
    "pflanzen",
    "loben",
    "graben",
    "filmen",
    "schlagen",
    "wirken",
    "klappen",
    "lenken",
    "wellen",
    "scheinen",
    "bremsen",
    "heben",
    "sichern",
    "decken",
    "passen",
    "gründen",
    "sinken",
    "streichen",
    "erben",
    "drehen",
    "tropfen",
    "schlafen",
    "ändern",
    "kichern",
    "kümmern",
    "wecken",
    "heißen",
    "nennen",
    "heilen",
    "gleichen",
    "tollen",
    "schwingen",
    "liegen",
    "formen",
    "quälen",
    "hauen",
    "doppeln",
    "mögen",
    "pfeifen",
    "pissen",
    "kehren",
    "zeigen",
    "lockern",
    "wenden",
    "löschen",
    "werden",
    "rauben",
    "reizen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

